$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.6946179502469849
$ws.Range("D2").Value = 0.1674370704476836
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "3.10e-13"

# Row 3
$ws.Range("C3").Value = 0.6520795768221956
$ws.Range("D3").Value = 0.1707295210231295
$ws.Range("E3").Value = "qa_saccade_regression_rate_%"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "1.61e-09"

# Row 4
$ws.Range("C4").Value = 0.7389585117874775
$ws.Range("D4").Value = 0.1395684625457408
$ws.Range("E4").Value = "qa_saccade_regression_rate_%"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "2.18e-07"
